$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D11").Value = "[1, 0, 1, 0, 1, 0, 0]"
$ws.Range("E11").Value = "['Normal', 'HardwareFault', 'RegulationViolation']"

$ws.Range("D12").Value = "[1, 0, 1, 0, 0, 0, 0]"
$ws.Range("E12").Value = "['Normal', 'HardwareFault']"

$ws.Range("D15").Value = "[0, 0, 0, 1, 0, 0, 0]"
$ws.Range("E15").Value = "['ParamViolation']"

$ws.Range("D25").Value = "[1, 0, 0, 0, 0, 0, 1]"
$ws.Range("E25").Value = "['Normal', 'SoftwareFault']"

$ws.Range("D27").Value = "[0, 0, 0, 0, 0, 0, 1]"
$ws.Range("E27").Value = "['SoftwareFault']"

$ws.Range("D31").Value = "[0, 0, 0, 0, 0, 0, 1]"
$ws.Range("E31").Value = "['SoftwareFault']"

$ws.Range("D38").Value = "[1, 0, 1, 0, 0, 0, 1]"
$ws.Range("E38").Value = "['Normal', 'HardwareFault', 'SoftwareFault']"

$ws.Range("D56").Value = "[0, 0, 1, 0, 0, 0, 0]"
$ws.Range("E56").Value = "['HardwareFault']"

$ws.Range("D58").Value = "[0, 0, 0, 1, 0, 0, 0]"
$ws.Range("E58").Value = "['ParamViolation']"

$ws.Range("D69").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E69").Value = "['Normal']"

$ws.Range("D73").Value = "[1, 0, 0, 1, 0, 0, 0]"
$ws.Range("E73").Value = "['Normal', 'ParamViolation']"

$ws.Range("D75").Value = "[0, 0, 0, 0, 0, 0, 0]"
$ws.Range("E75").Value = "[]"

$ws.Range("D83").Value = "[1, 1, 0, 0, 0, 0, 0]"
$ws.Range("E83").Value = "['Normal', 'SurroundingEnvironment']"

$ws.Range("D92").Value = "[1, 0, 1, 0, 0, 0, 1]"
$ws.Range("E92").Value = "['Normal', 'HardwareFault', 'SoftwareFault']"

$ws.Range("D97").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E97").Value = "['Normal']"

$ws.Range("D107").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E107").Value = "['Normal']"

$ws.Range("D109").Value = "[1, 1, 0, 0, 0, 0, 1]"
$ws.Range("E109").Value = "['Normal', 'SurroundingEnvironment', 'SoftwareFault']"

$ws.Range("D118").Value = "[1, 0, 0, 0, 1, 0, 0]"
$ws.Range("E118").Value = "['Normal', 'RegulationViolation']"
